$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp bumped from 20:16 to 20:46
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 20:46"

# Country column (A) is positional/rank-sorted, so several rows below
# both shift which country they hold AND refresh that countrys daily
# figures (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) in columns B..H.
$updates = @(
  @{ Row = 6; Country = "Estados Unidos"; Vals = @(51768, 8034, 370, 50731, 1175, 114, 667) },
  @{ Row = 22; Country = "Brasil"; Vals = @(2018, 94, 2, 1982, 18, 0, 34) },
  @{ Row = 23; Country = "Turquia"; Vals = @(1872, 343, 0, 1828, 0, 7, 44) },
  @{ Row = 24; Country = "Israel"; Vals = @(1656, 214, 49, 1604, 34, 2, 3) },
  @{ Row = 25; Country = "Malasia"; Vals = @(1624, 106, 183, 1425, 64, 2, 16) },
  @{ Row = 26; Country = "Dinamarca"; Vals = @(1591, 131, 1, 1558, 69, 8, 32) },
  @{ Row = 38; Country = "Rumania"; Vals = @(762, 186, 79, 672, 15, 4, 11) },
  @{ Row = 94; Country = "Reunion"; Vals = @(94, 23, 1, 93, 0, 0, 0) },
  @{ Row = 95; Country = "Camboya"; Vals = @(91, 4, 4, 87, 1, 0, 0) },
  @{ Row = 96; Country = "Azerbaiyan"; Vals = @(87, 15, 10, 76, 6, 0, 1) },
  @{ Row = 97; Country = "Senegal"; Vals = @(86, 7, 8, 78, 0, 0, 0) },
  @{ Row = 98; Country = "Venezuela"; Vals = @(84, 0, 15, 69, 2, 0, 0) },
  @{ Row = 99; Country = "Oman"; Vals = @(84, 18, 17, 67, 0, 0, 0) },
  @{ Row = 110; Country = "Ghana"; Vals = @(53, 26, 0, 51, 0, 0, 2) },
  @{ Row = 113; Country = "Cuba"; Vals = @(48, 8, 1, 46, 2, 0, 1) },
  @{ Row = 135; Country = "Guyana"; Vals = @(20, 0, 0, 19, 0, 0, 1) },
  @{ Row = 136; Country = "Togo"; Vals = @(20, 2, 1, 19, 0, 0, 0) },
  @{ Row = 137; Country = "Madagascar"; Vals = @(17, 5, 0, 17, 0, 0, 0) },
  @{ Row = 138; Country = "Islas Virgenes de los Estados Unidos"; Vals = @(17, 0, 0, 17, 0, 0, 0) },
  @{ Row = 142; Country = "Tanzania"; Vals = @(12, 0, 0, 12, 0, 0, 0) },
  @{ Row = 143; Country = "Etiopia"; Vals = @(12, 1, 0, 12, 0, 0, 0) },
  @{ Row = 147; Country = "Uganda"; Vals = @(9, 0, 0, 9, 0, 0, 0) },
  @{ Row = 148; Country = "Guinea Ecuatorial"; Vals = @(9, 0, 0, 9, 0, 0, 0) },
  @{ Row = 150; Country = "Surinam"; Vals = @(7, 2, 0, 7, 0, 0, 0) },
  @{ Row = 151; Country = "Seychelles"; Vals = @(7, 0, 0, 7, 0, 0, 0) },
  @{ Row = 152; Country = "Benin"; Vals = @(6, 0, 0, 6, 0, 0, 0) },
  @{ Row = 153; Country = "Bermudas"; Vals = @(6, 0, 0, 6, 0, 0, 0) },
  @{ Row = 154; Country = "Haiti"; Vals = @(6, 0, 0, 6, 0, 0, 0) },
  @{ Row = 160; Country = "Santa Sede"; Vals = @(4, 3, 0, 4, 0, 0, 0) },
  @{ Row = 161; Country = "Fiyi"; Vals = @(4, 1, 0, 4, 0, 0, 0) },
  @{ Row = 162; Country = "Congo"; Vals = @(4, 0, 0, 4, 0, 0, 0) },
  @{ Row = 163; Country = "Suazilandia"; Vals = @(4, 0, 0, 4, 0, 0, 0) },
  @{ Row = 164; Country = "Guinea"; Vals = @(4, 0, 0, 4, 0, 0, 0) },
  @{ Row = 165; Country = "Bahamas"; Vals = @(4, 0, 0, 4, 0, 0, 0) },
  @{ Row = 173; Country = "Birmania"; Vals = @(3, 1, 0, 3, 0, 0, 0) },
  @{ Row = 174; Country = "Santa Lucia"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
  @{ Row = 175; Country = "Zambia"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
  @{ Row = 176; Country = "Republica de Africa Central"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
  @{ Row = 177; Country = "Republica de Yibuti"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
  @{ Row = 178; Country = "San Bartolome"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
  @{ Row = 179; Country = "Sudan"; Vals = @(3, 1, 0, 2, 0, 0, 1) },
  @{ Row = 180; Country = "Gambia"; Vals = @(3, 1, 0, 2, 0, 0, 1) },
  @{ Row = 181; Country = "Cabo Verde"; Vals = @(3, 0, 0, 2, 0, 1, 1) },
  @{ Row = 182; Country = "Zimbabue"; Vals = @(3, 0, 0, 2, 0, 0, 1) },
  @{ Row = 183; Country = "Mauritania"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
  @{ Row = 184; Country = "Butan"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
  @{ Row = 185; Country = "Laos"; Vals = @(2, 2, 0, 2, 0, 0, 0) },
  @{ Row = 186; Country = "San Martin (Parte Holandesa)"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
  @{ Row = 187; Country = "Dominica"; Vals = @(2, 0, 0, 2, 0, 0, 0) }
)

$cols = @("B","C","D","E","F","G","H")
foreach ($item in $updates) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Country
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $item.Vals[$i]
    }
}
